$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.448.36"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.629.98"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "304.20"
$ws.Range("E6").Value = "  -0.87%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3771"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3655"
$ws.Range("E8").Value = "  +0.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "51.59"
$ws.Range("E9").Value = "  -1.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08249"
$ws.Range("E10").Value = "  +0.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.228"
$ws.Range("E11").Value = "  -3.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.45"
$ws.Range("E13").Value = "  -2.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.564"
$ws.Range("E14").Value = "  -1.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001253"
$ws.Range("E15").Value = "  -2.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.309"
$ws.Range("E16").Value = "  -0.99%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.629.54"
$ws.Range("E17").Value = "  -0.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.07"
$ws.Range("E18").Value = "  -0.74%  "
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.74"
$ws.Range("E20").Value = "  -2.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.479"
$ws.Range("E21").Value = "  -1.01%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.72"
$ws.Range("E23").Value = "  -0.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.437.72"
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.196"
$ws.Range("E25").Value = "  +3.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.474"
$ws.Range("E26").Value = "  +2.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.40"
$ws.Range("E27").Value = "  +0.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "149.82"
$ws.Range("E28").Value = "  -0.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.318"
$ws.Range("E29").Value = "  -0.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.97"
$ws.Range("E30").Value = "  -1.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.811.21"
$ws.Range("E31").Value = "  -0.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.273"
$ws.Range("E32").Value = "  -3.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.802"
$ws.Range("E33").Value = "  +0.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.019"
$ws.Range("E34").Value = "  +5.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "10.85"
$ws.Range("E35").Value = "  +4.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02787"
$ws.Range("E36").Value = "  -1.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2526"
$ws.Range("E37").Value = "  -0.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.08751"
$ws.Range("E38").Value = "  -1.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.07136"
$ws.Range("E39").Value = "  -2.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.037"
$ws.Range("E40").Value = "  -2.35%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.354"
$ws.Range("E41").Value = "  -1.90%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7054"
$ws.Range("E42").Value = "  -0.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.26"
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.23"
$ws.Range("E44").Value = "  -2.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6561"
$ws.Range("E45").Value = "  +0.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.325"
$ws.Range("E46").Value = "  -0.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9999"
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.979"
$ws.Range("E48").Value = "  -1.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08019"
$ws.Range("E49").Value = "  +0.57%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.200"
$ws.Range("E50").Value = "  -0.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "126.20"
$ws.Range("E51").Value = "  -2.68%  "
